$d = $word.ActiveDocument

# --- 1. Change the language tags on the run that holds the inline picture ---
# <w:lang w:eastAsia="da-DK"/>  ->  <w:lang w:val="en-GB" w:eastAsia="en-GB"/>
$picRange = $d.InlineShapes(1).Range
$picRange.LanguageID = 2057        # wdEnglishUK
$picRange.LanguageIDFarEast = 2057 # wdEnglishUK

# --- 2. Remove ", fjerne" from the recipe paragraph ---
$d.Content.Find.Execute(
    "samt ændre, fjerne og finde allerede eksisterende opskrifter",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "samt ændre og finde allerede eksisterende opskrifter", 2) | Out-Null

# --- 3. Append a new sentence about the administrator removing recipes ---
$target = "også oprette en madplan for en uge med opskrifter og generere en indkøbsliste ud fra det."
$found = $d.Content.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$endRange = $d.Content
$endRange.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $endRange.Duplicate
$insertPoint.Collapse(0)  # wdCollapseEnd
$insertPoint.InsertAfter(" I Pristjek220 Forretning er det administratoren, der kan fjerne opskrifter hvis dette ønskes.")

# --- 4. Remove the _GoBack bookmark from the end of the document, it now lives after the new sentence ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-add the _GoBack bookmark right after the newly inserted sentence
$newEnd = $insertPoint.Duplicate
$newEnd.Collapse(0)
$newEnd.MoveEnd(1, 0)
$d.Bookmarks.Add("_GoBack", $newEnd)
